$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the style used by the
# existing header row (B1:H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New data columns I2:J36.
$iVals = @(9,4,7,6,8,8,8,6,7,6,7,6,7,6,5,9,7,10,8,4,7,9,7,3,7,4,6,8,1,9,9,7,3,7,5)
$jVals = @(9,6,7,6,8,8,8,6,7,7,7,7,7,6,6,9,7,10,8,5,7,9,8,4,8,4,6,8,1,9,9,7,4,7,5)

for ($r = 2; $r -le 36; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iVals[$idx]
    $ws.Cells.Item($r, 10).Value = $jVals[$idx]
}

Write-Output "done"
